# Update "想去人数" (number of people interested) counts for several
# events on the "展览" sheet and the matching rows on the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row => new value, for the "展览" sheet (column F)
$exhibitUpdates = @{
    3  = 547
    11 = 67
    14 = 511
    15 = 22
    16 = 6527
    22 = 15586
    24 = 297
    27 = 11121
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => new value, for the "全部类型" sheet (column F)
$allUpdates = @{
    3  = 547
    13 = 67
    16 = 511
    18 = 22
    19 = 6527
    26 = 15586
    28 = 297
    32 = 11121
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
